$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: find a paragraph index whose Range.Text starts with a given
# marker string (used instead of hard-coded indices so the script is
# resilient to earlier insertions shifting later paragraph numbers).
# ---------------------------------------------------------------------
function Find-ParagraphIndex($marker) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs($i).Range.Text
        if ($t -like "$marker*") {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# 1) Insert the new "Note:" paragraph plus the two "sudo" command
#    paragraphs right after the "PC Users: In future assignment
#    walkthroughs..." paragraph (and before "Accessing User Interfaces").
# ---------------------------------------------------------------------
$anchorIdx = Find-ParagraphIndex("PC Users: In future assignment walkthroughs")
$anchorRng = $d.Paragraphs($anchorIdx).Range
$anchorRng.InsertParagraphAfter()

$insertIdx = $anchorIdx + 1
$insertRng = $d.Paragraphs($insertIdx).Range

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newParasXml = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="BodyText"/>
    <w:rPr>
      <w:color w:val="EE0000"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:color w:val="EE0000"/>
    </w:rPr>
    <w:t>Note:</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:color w:val="EE0000"/>
    </w:rPr>
    <w:t xml:space="preserve"> Over time, Ubuntu releases will change. Always start with the latest minimal install. Some minimal images do not include Git by default, so you’ll need to install it manually with:</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="BodyText"/>
    <w:rPr>
      <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
      <w:color w:val="EE0000"/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
      <w:color w:val="EE0000"/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t>sudo</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
      <w:color w:val="EE0000"/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t xml:space="preserve"> apt update</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="BodyText"/>
    <w:rPr>
      <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
      <w:color w:val="EE0000"/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
      <w:color w:val="EE0000"/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t>sudo</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
      <w:color w:val="EE0000"/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t xml:space="preserve"> apt install git -y</w:t>
  </w:r>
</w:p>
"@

$insertRng.InsertXML($newParasXml)

# The "Note:" run needs a character style of "Strong" in addition to the
# red color override - InsertXML drops rStyle, so re-apply it afterwards
# via Range.Style on the precise "Note:" sub-range.
$noteParaIdx = $insertIdx
$noteRng = $d.Paragraphs($noteParaIdx).Range
$noteStrongRng = $d.Range($noteRng.Start, $noteRng.Start + 5)
$noteStrongRng.Style = "Strong"

# ---------------------------------------------------------------------
# 2) Add <w:lastRenderedPageBreak/> as the first child of the run that
#    holds "Stop " in the "Stop NiFi:" bullet (right before ./nifi.sh
#    stop). We rebuild the whole paragraph (preserving its original
#    identifiers) via InsertXML on its own Range, which replaces it in
#    place without disturbing paragraph count.
# ---------------------------------------------------------------------
$stopIdx = Find-ParagraphIndex("Stop NiFi")
$stopRng = $d.Paragraphs($stopIdx).Range

$stopParaXml = @"
<w:p $wns xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="648D4B54" w14:textId="5B37CE4B" w:rsidR="00944F9A" w:rsidRPr="00DE2081" w:rsidRDefault="008E7D89" w:rsidP="00944F9A">
  <w:pPr>
    <w:pStyle w:val="Compact"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="22"/>
    </w:numPr>
  </w:pPr>
  <w:r w:rsidRPr="008E7D89">
    <w:lastRenderedPageBreak/>
    <w:t xml:space="preserve">Stop </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r w:rsidRPr="008E7D89">
    <w:t>NiFi</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r w:rsidRPr="008E7D89">
    <w:t>:</w:t>
  </w:r>
  <w:r w:rsidRPr="008E7D89">
    <w:br/>
  </w:r>
  <w:r w:rsidRPr="008E7D89">
    <w:rPr>
      <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
    </w:rPr>
    <w:t>./nifi.sh stop</w:t>
  </w:r>
  <w:bookmarkEnd w:id="0"/>
  <w:bookmarkEnd w:id="5"/>
</w:p>
"@

$stopRng.InsertXML($stopParaXml)

Write-Host "Edit complete"
